$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row 16 (shifts Spain..South Africa down by one row)
$ws.Rows.Item(16).Insert()

# Populate the new row with Egypt's data
$ws.Range("A16").Value = "Egypt"
$ws.Range("B16").Value = "EGY"
$ws.Range("C16").Value = "EG"
$ws.Range("D16").Value = 30
$ws.Range("E16").Value = 31.2357
$ws.Range("F16").Value = 30.0444

# Match the longitude/latitude cell formatting used by other manually-added
# rows (plain style, not the 3-decimal number format) by copying the format
# from the Corruption Index cell in the same row.
$ws.Range("D16").Copy()
$ws.Range("E16:F16").PasteSpecial(-4122)

# Select cell G1, as in the saved workbook after the edit
$ws.Range("G1").Select()
